# "Fix outdated help command"
#
# 1) The deck's auto-updating "date" placeholders (slide master, every
#    custom layout, and the notes master) were re-cached by PowerPoint from
#    15/11/2021 to 8/01/2022 the next time the file was saved.
# 2) Slide 15's "singularity help" code sample was corrected to
#    "singularity run-help".

$p = $ppt.ActivePresentation
$newDate = "8/01/2022"
$ppPlaceholderDate = 16

function Set-DateOnShapes($shapes, $newText) {
    if ($shapes -eq $null) { return }
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Set-DateOnShapes $master.Shapes $newDate

# Every custom (slide) layout hanging off the master.
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Set-DateOnShapes $layout.Shapes $newDate
}

# Notes master's date placeholder.
Set-DateOnShapes $p.NotesMaster.Shapes $newDate

# Slide 15 ("Recipe file: tests & meta-info"): fix the stale help command,
# leaving every other run / its Courier New formatting untouched.
$slide = $p.Slides.Item(15)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$fullText = $tr.Text
$needle = "singularity help"
$idx = $fullText.IndexOf($needle)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $needle.Length)
    $sub.Text = "singularity run-help"
}
